$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A160").Value = "IMX-USD"
$ws.Range("A161").Value = "GRT-USD"
